$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row from the old Korean/mixed labels to the new
# lowercase field-name headers (data read error handle: align headers
# with the JSON field names the app actually reads).
$ws.Range("A1").Value = "korName"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "weight"
$ws.Range("D1").Value = "volume"
$ws.Range("E1").Value = "description"

# Move the active selection to C4, matching the author's last cursor spot.
$ws.Range("C4").Select()
